# "Merge Code from Master to NoCodeWait"
# Updates the Config sheet's RunInParallel flag to "Yes" and the
# "Execute" note on the first test case row, then restores the
# selection/active-sheet state captured in the target workbook.

$wb = $excel.ActiveWorkbook

$wsConfig = $wb.Worksheets.Item("Config")
$wsTest   = $wb.Worksheets.Item("Test Cases")

# --- Data changes -----------------------------------------------------

# Config!B3 (RunInParallel) : "No" -> "Yes"
$wsConfig.Range("B3").Value = "Yes"

# Test Cases!D2 (Execute note for TestCaseNumber 101)
# "TestCaseNumber=101" -> "Groups=Regression"
$wsTest.Range("D2").Value = "Groups=Regression"

# --- Selection / active sheet changes ----------------------------------

# "Test Cases" is no longer the active tab; its remembered selection moves to D7
$wsTest.Activate() | Out-Null
$wsTest.Range("D7").Select() | Out-Null

# "Config" becomes the active tab, with B4 selected
$wsConfig.Activate() | Out-Null
$wsConfig.Range("B4").Select() | Out-Null
